# Insert a new weekly record for "Feria Lagunitas de Puerto Montt - Betarraga"
# at row 178, pushing the existing rows 178-211 down to 179-212.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(178).Insert()

$ws.Range("A178").Value = 4
$ws.Range("B178").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value = "Los Lagos"
$ws.Range("D178").Value = 44522
$ws.Range("E178").Value = 10
$ws.Range("F178").Value = 100114014
$ws.Range("G178").Value = "Betarraga"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 500
$ws.Range("K178").Value = 1000
$ws.Range("L178").Value = 1000
$ws.Range("M178").Value = 1000
$ws.Range("N178").Value = '$/paquete 5 unidades'
$ws.Range("O178").Value = "Región del Maule"
$ws.Range("P178").Value = 200
$ws.Range("Q178").Value = 5
$ws.Range("R178").Value = "Hortaliza"
